{"js": "// Update the date title and the twenty-five \"three-digit \u00f7 one-digit\"\n// practice answers laid out in the 5-column table beneath it.\n//\n// The replacements are applied positionally (title first, then each\n// table cell in row-major / left-to-right order, skipping the blank\n// spacer rows) because one answer string (\"880\u00f73=293, 1\") appears twice\n// in the source but maps to two different replacements depending on\n// where it sits in the table \u2014 a global text find/replace would not be\n// able to tell those two occurrences apart.\nconst replacements = [\n  \"2026-01-13 Tuesday\",\n  \"430\u00f78=53, 6\",\n  \"904\u00f78=113, 0\",\n  \"527\u00f72=263, 1\",\n  \"794\u00f78=99, 2\",\n  \"383\u00f78=47, 7\",\n  \"846\u00f72=423, 0\",\n  \"717\u00f72=358, 1\",\n  \"598\u00f79=66, 4\",\n  \"514\u00f78=64, 2\",\n  \"476\u00f72=238, 0\",\n  \"564\u00f72=282, 0\",\n  \"625\u00f73=208, 1\",\n  \"762\u00f77=108, 6\",\n  \"390\u00f78=48, 6\",\n  \"696\u00f74=174, 0\",\n  \"814\u00f75=162, 4\",\n  \"333\u00f77=47, 4\",\n  \"192\u00f72=96, 0\",\n  \"204\u00f74=51, 0\",\n  \"856\u00f77=122, 2\",\n  \"542\u00f74=135, 2\",\n  \"813\u00f73=271, 0\",\n  \"187\u00f72=93, 1\",\n  \"478\u00f79=53, 1\",\n  \"449\u00f73=149, 2\"\n];\nlet next = 0;\n\nconst body = context.document.body;\n\n// 1) The title paragraph (\"2026-01-12 Monday\" -> \"2026-01-13 Tuesday\")\n// is the very first paragraph in the document body.\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\nparagraphs.items[0].insertText(replacements[next++], \"Replace\");\n\n// 2) The answer grid is the single table in the document body.\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\nfor (const row of rows.items) {\n  const cells = row.cells;\n  cells.load(\"items\");\n  await context.sync();\n\n  for (const cell of cells.items) {\n    const range = cell.getRange();\n    range.load(\"text\");\n    await context.sync();\n\n    if (range.text.trim().length === 0) continue;\n    if (next >= replacements.length) continue;\n\n    cell.value = replacements[next++];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date title and the twenty-five \"three-digit \u00f7 one-digit\"\n# practice answers laid out in the 5-column table beneath it.\n#\n# The replacements are applied positionally (title first, then each\n# table cell in row-major / left-to-right order, skipping the blank\n# spacer rows) because one answer string (\"880\u00f73=293, 1\") appears twice\n# in the source but maps to two different replacements depending on\n# where it sits in the table - a global text find/replace would not be\n# able to tell those two occurrences apart.\n\n$replacements = @(\n    \"2026-01-13 Tuesday\",\n    \"430\u00f78=53, 6\",\n    \"904\u00f78=113, 0\",\n    \"527\u00f72=263, 1\",\n    \"794\u00f78=99, 2\",\n    \"383\u00f78=47, 7\",\n    \"846\u00f72=423, 0\",\n    \"717\u00f72=358, 1\",\n    \"598\u00f79=66, 4\",\n    \"514\u00f78=64, 2\",\n    \"476\u00f72=238, 0\",\n    \"564\u00f72=282, 0\",\n    \"625\u00f73=208, 1\",\n    \"762\u00f77=108, 6\",\n    \"390\u00f78=48, 6\",\n    \"696\u00f74=174, 0\",\n    \"814\u00f75=162, 4\",\n    \"333\u00f77=47, 4\",\n    \"192\u00f72=96, 0\",\n    \"204\u00f74=51, 0\",\n    \"856\u00f77=122, 2\",\n    \"542\u00f74=135, 2\",\n    \"813\u00f73=271, 0\",\n    \"187\u00f72=93, 1\",\n    \"478\u00f79=53, 1\",\n    \"449\u00f73=149, 2\"\n)\n$next = 0\n\n$d = $word.ActiveDocument\n\n# 1) The title paragraph (\"2026-01-12 Monday\" -> \"2026-01-13 Tuesday\")\n# is the very first paragraph in the document body.\n$titlePara = $d.Paragraphs.Item(1)\n$titlePara.Range.Text = $replacements[$next]\n$next = $next + 1\n\n# 2) The answer grid is the single table in the document body.\n$t = $d.Tables.Item(1)\n$rowCount = $t.Rows.Count\n$colCount = $t.Columns.Count\n\nfor ($r = 1; $r -le $rowCount; $r++) {\n    for ($c = 1; $c -le $colCount; $c++) {\n        $cell = $t.Cell($r, $c)\n        $raw = $cell.Range.Text\n        # Cell.Range.Text always carries a trailing cell-mark (CR + BEL);\n        # strip it before checking whether the cell actually has content.\n        $clean = $raw.TrimEnd([char]7).TrimEnd([char]13)\n        if ($clean.Length -eq 0) {\n            continue\n        }\n        if ($next -ge $replacements.Length) {\n            continue\n        }\n        $cell.Range.Text = $replacements[$next]\n        $next = $next + 1\n    }\n}\n"}
